$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1272.7778
$ws.Range("I19").Value = 1134
$ws.Range("K19").Value = 1134
$ws.Range("M19").Value = -959
$ws.Range("H32").Value = 1573.9375
$ws.Range("I32").Value = 1570.5714
$ws.Range("K32").Value = 1570.5714
$ws.Range("M32").Value = -1244.5714
$ws.Range("H64").Value = 7885.467
$ws.Range("I64").Value = 4328.3
$ws.Range("K64").Value = 4328.3
$ws.Range("M64").Value = -4080.3
$ws.Range("H67").Value = 7885.467
$ws.Range("I67").Value = 4328.3
$ws.Range("K67").Value = 4328.3
$ws.Range("M67").Value = -3470.3
$ws.Range("H70").Value = 25855
$ws.Range("I70").Value = 1700
$ws.Range("K70").Value = 5100
$ws.Range("M70").Value = -4830
$ws.Range("H73").Value = 25855
$ws.Range("I73").Value = 1700
$ws.Range("K73").Value = 5100
$ws.Range("M73").Value = -4164
$ws.Range("H74").Value = 2999.889
$ws.Range("I74").Value = 2999.5
$ws.Range("K74").Value = 2999.5
$ws.Range("M74").Value = -2063.5
$ws.Range("H77").Value = 2999.889
$ws.Range("I77").Value = 2999.5
$ws.Range("K77").Value = 14997.5
$ws.Range("M77").Value = -10317.5
$ws.Range("H98").Value = 1267.8572
$ws.Range("I98").Value = 1267.8572
$ws.Range("K98").Value = 1267.8572
$ws.Range("M98").Value = 230.1428000000001
$ws.Range("H122").Value = 1267.8572
$ws.Range("I122").Value = 1267.8572
$ws.Range("K122").Value = 3803.5716
$ws.Range("M122").Value = -1353.5716

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3325.9375
$ws.Range("I32").Value = 2360.1155
$ws.Range("K32").Value = 2360.1155
$ws.Range("M32").Value = -2073.1155
$ws.Range("H45").Value = 4259.4165
$ws.Range("I45").Value = 1778.25
$ws.Range("J45").Value = 5500
$ws.Range("K45").Value = 1778.25
$ws.Range("L45").Value = 5500
$ws.Range("M45").Value = -1401.25
$ws.Range("N45").Value = -6254
$ws.Range("H132").Value = 1424.6666
$ws.Range("I132").Value = 944.4286
$ws.Range("K132").Value = 2833.2858
$ws.Range("M132").Value = -303.2857999999997

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2686.8823
$ws.Range("I105").Value = 2611.8667
$ws.Range("J105").Value = 3249.5
$ws.Range("K105").Value = 2611.8667
$ws.Range("L105").Value = 3249.5
$ws.Range("M105").Value = -864.8667
$ws.Range("N105").Value = -6743.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1916.5714
$ws.Range("I31").Value = 1983.2
$ws.Range("K31").Value = 1983.2
$ws.Range("M31").Value = -1688.2
$ws.Range("H34").Value = 1916.5714
$ws.Range("I34").Value = 1983.2
$ws.Range("K34").Value = 1983.2
$ws.Range("M34").Value = -1781.2
$ws.Range("H41").Value = 5474.6665
$ws.Range("J41").Value = 9975
$ws.Range("L41").Value = 9975
$ws.Range("N41").Value = -10831
$ws.Range("H50").Value = 20041.5
$ws.Range("J50").Value = 20000
$ws.Range("L50").Value = 20000
$ws.Range("N50").Value = -21250
$ws.Range("H58").Value = 1050.091
$ws.Range("I58").Value = 672.3889
$ws.Range("K58").Value = 672.3889
$ws.Range("M58").Value = -469.3889
$ws.Range("H62").Value = 2000
$ws.Range("I62").Value = 2000
$ws.Range("K62").Value = 2000
$ws.Range("M62").Value = -1376
$ws.Range("H65").Value = 2000
$ws.Range("I65").Value = 2000
$ws.Range("K65").Value = 10000
$ws.Range("M65").Value = -6880
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H74").Value = 35000
$ws.Range("J74").Value = 35000
$ws.Range("L74").Value = 35000
$ws.Range("N74").Value = -36748
$ws.Range("H77").Value = 35000
$ws.Range("J77").Value = 35000
$ws.Range("L77").Value = 105000
$ws.Range("N77").Value = -113736
$ws.Range("H99").Value = 7207.1
$ws.Range("I99").Value = 7063.4443
$ws.Range("K99").Value = 7063.4443
$ws.Range("M99").Value = -5565.4443
$ws.Range("H107").Value = 561.625
$ws.Range("I107").Value = 499
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 499
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 1421
$ws.Range("N107").Value = -4840
$ws.Range("H126").Value = 7207.1
$ws.Range("I126").Value = 7063.4443
$ws.Range("K126").Value = 21190.3329
$ws.Range("M126").Value = -18720.3329
$ws.Range("H136").Value = 1050.091
$ws.Range("I136").Value = 672.3889
$ws.Range("K136").Value = 2017.1667
$ws.Range("M136").Value = 532.8332999999998

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 272.8889
$ws.Range("I7").Value = 37.333332
$ws.Range("J7").Value = 390.66666
$ws.Range("K7").Value = 111.999996
$ws.Range("L7").Value = 1171.99998
$ws.Range("M7").Value = 0.000004000000004111826
$ws.Range("N7").Value = -1395.99998
$ws.Range("H12").Value = 302.42856
$ws.Range("J12").Value = 302.42856
$ws.Range("L12").Value = 907.28568
$ws.Range("N12").Value = -1253.28568
$ws.Range("H113").Value = 842.4286
$ws.Range("I113").Value = 761.5
$ws.Range("K113").Value = 2284.5
$ws.Range("M113").Value = -114.5
$ws.Range("H122").Value = 582.6667
$ws.Range("I122").Value = 374.5
$ws.Range("K122").Value = 3370.5
$ws.Range("M122").Value = -920.5
$ws.Range("H140").Value = 1343.3334
$ws.Range("I140").Value = 1343.3334
$ws.Range("K140").Value = 4030.0002
$ws.Range("M140").Value = 1149.9998

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3573
$ws.Range("I126").Value = 3003.6667
$ws.Range("K126").Value = 9011.000100000001
$ws.Range("M126").Value = -6541.000100000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1480.5
$ws.Range("I16").Value = 1500.75
$ws.Range("J16").Value = 1399.5
$ws.Range("K16").Value = 1500.75
$ws.Range("L16").Value = 1399.5
$ws.Range("M16").Value = -1330.75
$ws.Range("N16").Value = -1739.5
$ws.Range("H74").Value = 34999.668
$ws.Range("J74").Value = 47499.5
$ws.Range("L74").Value = 47499.5
$ws.Range("N74").Value = -49495.5
$ws.Range("H77").Value = 34999.668
$ws.Range("J77").Value = 47499.5
$ws.Range("L77").Value = 142498.5
$ws.Range("N77").Value = -152482.5
$ws.Range("H80").Value = 25000
$ws.Range("J80").Value = 25000
$ws.Range("L80").Value = 25000
$ws.Range("N80").Value = -27246
$ws.Range("H83").Value = 25000
$ws.Range("J83").Value = 25000
$ws.Range("L83").Value = 75000
$ws.Range("N83").Value = -86232

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 20000
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 20000
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 20000
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -20780
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H100").Value = 5809623.5
$ws.Range("I100").Value = 7744776.5
$ws.Range("K100").Value = 15489553
$ws.Range("M100").Value = -15489012
$ws.Range("H113").Value = 525.5
$ws.Range("I113").Value = 509
$ws.Range("K113").Value = 1527
$ws.Range("M113").Value = 643
